$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034462150653301
$ws.Range("D2").Value = 1.041505736932444
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.048827187618355
$ws.Range("I2").Value = 1.033845067088311
$ws.Range("J2").Value = 1.039581263592506
$ws.Range("K2").Value = 1.044285060598727
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.051585947791
$ws.Range("N2").Value = 1.017070309472182
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035706698738866
$ws.Range("D3").Value = 1.042481051107273
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.050008122977228
$ws.Range("I3").Value = 1.034078485836317
$ws.Range("J3").Value = 1.040467461277363
$ws.Range("K3").Value = 1.045070583984231
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.052578053233864
$ws.Range("N3").Value = 1.017370826705859
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036511427353883
$ws.Range("D4").Value = 1.04311145718692
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.05077199366129
$ws.Range("I4").Value = 1.034227867082601
$ws.Range("J4").Value = 1.041039836076952
$ws.Range("K4").Value = 1.045577569774256
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.05321916666639
$ws.Range("N4").Value = 1.017564752689465
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036849599859392
$ws.Range("D5").Value = 1.043376316642288
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.051093061258826
$ws.Range("I5").Value = 1.034290270839701
$ws.Range("J5").Value = 1.041280211681997
$ws.Range("K5").Value = 1.045790396577484
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.053488490043443
$ws.Range("N5").Value = 1.017646153058475
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036906372651667
$ws.Range("D6").Value = 1.043420778174416
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.051146966242459
$ws.Range("I6").Value = 1.034300725488092
$ws.Range("J6").Value = 1.041320557147029
$ws.Range("K6").Value = 1.04582611299161
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.053533698900608
$ws.Range("N6").Value = 1.017659813156853
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036515946560035
$ws.Range("D7").Value = 1.043114996891816
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.050776284030027
$ws.Range("I7").Value = 1.034228702480642
$ws.Range("J7").Value = 1.04104304897082
$ws.Range("K7").Value = 1.045580414793008
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.053222766165191
$ws.Range("N7").Value = 1.017565840860139
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034882872699682
$ws.Range("D8").Value = 1.041835492046228
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.049226347977285
$ws.Range("I8").Value = 1.033924295235742
$ws.Range("J8").Value = 1.039880977560616
$ws.Range("K8").Value = 1.044550801931867
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.051921410435243
$ws.Range("N8").Value = 1.017171980447228
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032000642045944
$ws.Range("D9").Value = 1.039575517007541
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.046492985655604
$ws.Range("I9").Value = 1.03337518724112
$ws.Range("J9").Value = 1.037825109389165
$ws.Range("K9").Value = 1.042726474935679
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.049621709658737
$ws.Range("N9").Value = 1.016473876964819
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030075921333623
$ws.Range("D10").Value = 1.038065201533079
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.044669162017402
$ws.Range("I10").Value = 1.033000547126919
$ws.Range("J10").Value = 1.036448944187035
$ws.Range("K10").Value = 1.04150343757711
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.048084069978331
$ws.Range("N10").Value = 1.016005706450368
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029241688876423
$ws.Range("D11").Value = 1.037410328831817
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.043879023734325
$ws.Range("I11").Value = 1.032836285136052
$ws.Range("J11").Value = 1.035851699892765
$ws.Range("K11").Value = 1.040972211925587
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.04741716210829
$ws.Range("N11").Value = 1.015802319645746
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028931691691071
$ws.Range("D12").Value = 1.037166943482172
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.04358546704796
$ws.Range("I12").Value = 1.032774963758257
$ws.Range("J12").Value = 1.035629650880648
$ws.Range("K12").Value = 1.040774642528489
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.047169275416257
$ws.Range("N12").Value = 1.015726672025296
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028998192903538
$ws.Range("D13").Value = 1.037219156671057
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.043648438896093
$ws.Range("I13").Value = 1.032788131311646
$ws.Range("J13").Value = 1.035677290500093
$ws.Range("K13").Value = 1.040817033110453
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.047222455564382
$ws.Range("N13").Value = 1.015742903264526
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029216067005553
$ws.Range("D14").Value = 1.037390213317511
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.043854759569902
$ws.Range("I14").Value = 1.032831222561428
$ws.Range("J14").Value = 1.035833349465074
$ws.Range("K14").Value = 1.040955885858042
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.047396675144937
$ws.Range("N14").Value = 1.015796068648061
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029350289613495
$ws.Range("D15").Value = 1.037495588824187
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.043981871972607
$ws.Range("I15").Value = 1.032857731782181
$ws.Range("J15").Value = 1.035929475195289
$ws.Range("K15").Value = 1.041041404631355
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.047503995328592
$ws.Range("N15").Value = 1.015828812231138
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030131269105331
$ws.Range("D16").Value = 1.038108644240562
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.044721592032454
$ws.Range("I16").Value = 1.033011405613541
$ws.Range("J16").Value = 1.036488552572619
$ws.Range("K16").Value = 1.04153865851843
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.048128307100236
$ws.Range("N16").Value = 1.016019190476807
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030620936016587
$ws.Range("D17").Value = 1.038492956493503
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.045185487219474
$ws.Range("I17").Value = 1.033107254451268
$ws.Range("J17").Value = 1.036838882527672
$ws.Range("K17").Value = 1.041850131769691
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.048519625405161
$ws.Range("N17").Value = 1.016138430977545
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030906471981669
$ws.Range("D18").Value = 1.038717032990532
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.045456029983404
$ws.Range("I18").Value = 1.03316296463526
$ws.Range("J18").Value = 1.037043093410974
$ws.Range("K18").Value = 1.042031650501379
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.048747768792552
$ws.Range("N18").Value = 1.016207917741748
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031003819165695
$ws.Range("D19").Value = 1.038793422695819
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.045548271456651
$ws.Range("I19").Value = 1.033181927002528
$ws.Range("J19").Value = 1.037112701950364
$ws.Range("K19").Value = 1.042093516893202
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.048825541885015
$ws.Range("N19").Value = 1.016231600062651
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030568407533689
$ws.Range("D20").Value = 1.038451732370472
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.045135719755616
$ws.Range("I20").Value = 1.03309699114112
$ws.Range("J20").Value = 1.036801308969954
$ws.Range("K20").Value = 1.041816730029853
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.048477651631198
$ws.Range("N20").Value = 1.016125644246198
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029151912022694
$ws.Range("D21").Value = 1.037339845180665
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.043794005064453
$ws.Range("I21").Value = 1.032818541746537
$ws.Range("J21").Value = 1.035787399690934
$ws.Range("K21").Value = 1.040915004059929
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.04734537645991
$ws.Range("N21").Value = 1.015780415553043
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028260573687387
$ws.Range("D22").Value = 1.036639967080143
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.042950043422901
$ws.Range("I22").Value = 1.03264169211285
$ws.Range("J22").Value = 1.035148722370329
$ws.Range("K22").Value = 1.040346614300648
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.046632500413361
$ws.Range("N22").Value = 1.015562773588567
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028733159522064
$ws.Range("D23").Value = 1.037011061249002
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.043397479576491
$ws.Range("I23").Value = 1.032735612139675
$ws.Range("J23").Value = 1.035487411061484
$ws.Range("K23").Value = 1.040648065451339
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.047010502145111
$ws.Range("N23").Value = 1.015678205194744
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030592143136371
$ws.Range("D24").Value = 1.038470360044005
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.045158207652219
$ws.Range("I24").Value = 1.03310162929736
$ws.Range("J24").Value = 1.03681828724702
$ws.Range("K24").Value = 1.041831823327592
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.048496618100308
$ws.Range("N24").Value = 1.016131422218077
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03274632466461
$ws.Range("D25").Value = 1.040160414445131
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.04719989546121
$ws.Range("I25").Value = 1.033518653014296
$ws.Range("J25").Value = 1.038357577357198
$ws.Range("K25").Value = 1.043199302946184
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.050217023885289
$ws.Range("N25").Value = 1.016654839071033
